# Apply the CloudGraph architecture overview edit: insert several new
# paragraphs of content describing HBase/Hadoop best practices and the
# HGraph API, right before the "Imagine taking an average sized..."
# paragraph, and leave the rest of that sentence (and everything after)
# untouched.

$d = $word.ActiveDocument

# Locate the target paragraph by its distinctive leading text rather than
# a hard-coded index, since paragraph numbering can vary.
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "Imagine taking an average sized relational database*") {
        $target = $cand
        break
    }
}

$para1 = "Numerous best practices have evolved out of the HBase and Hadoop open-source software ecosystem. Several are quite restrictive such as the support for ACID transactions only across a single HBase row. Other critical best practices involve the use of column families and in particular the format and length of composite row and column keys, as these all can effect the general performance and especially the even distribution of data across regions in an HBase cluster."
$para2 = "The HGraph implementation encapsulates many HBase best practices in each of these areas and provides the user with an intuitive convenient and standards-based API generated from a user-provided, "
$para3 = "domain-specific business model. Complexities of terse and efficient physical row and column key generation are completely hidden and the client user is provided with a meaningful view of his/her business entities and attributes."

$block = "`r" + $para1 + "`r`r" + $para2 + "`r" + $para3 + "`r"

$insertionPoint = $target.Range.Duplicate
$insertionPoint.Collapse(1)
$insertionPoint.InsertBefore($block)

Write-Output "inserted"
